# Refresh the cryptos list (prices + 1h volume deltas), as produced by the
# scheduled "Updated cryptos list" GitHub Actions job. All Price/Volume(1h)
# cells are plain text in this sheet, so numeric-looking values are written
# with a leading apostrophe (classic "store as text") and then have their
# style reset to Normal so Excel's auto quote-prefix formatting doesn't
# leave a stray style on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.027.05'
$ws.Range("E2").Value = '  -3.94%  '
$ws.Range("D3").Value = '2.968.08'
$ws.Range("E3").Value = '  -3.68%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = "'542.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.07%  '
$ws.Range("D6").Value = "'129.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.58%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '2.969.46'
$ws.Range("E8").Value = '  -3.45%  '
$ws.Range("D9").Value = "'0.489"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.71%  '
$ws.Range("D10").Value = "'0.143"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.84%  '
$ws.Range("D11").Value = "'5.86"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -7.60%  '
$ws.Range("D12").Value = "'0.438"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.12%  '
$ws.Range("D13").Value = "'0.0000215"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.82%  '
$ws.Range("D14").Value = "'33.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.23%  '
$ws.Range("D15").Value = '3.458.79'
$ws.Range("E15").Value = '  -3.49%  '
$ws.Range("D16").Value = '61.158.85'
$ws.Range("E16").Value = '  -3.80%  '
$ws.Range("E17").Value = '  -2.70%  '
$ws.Range("D18").Value = '2.974.73'
$ws.Range("E18").Value = '  -3.49%  '
$ws.Range("D19").Value = "'6.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.49%  '
$ws.Range("D20").Value = "'467.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.42%  '
$ws.Range("D21").Value = "'12.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.19%  '
$ws.Range("D22").Value = "'0.654"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.76%  '
$ws.Range("D23").Value = "'6.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.42%  '
$ws.Range("D24").Value = "'79.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.27%  '
$ws.Range("D25").Value = "'11.79"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.06%  '
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").Value = "'2.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.65%  '
$ws.Range("D28").Value = "'7.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.35%  '
$ws.Range("E29").Value = '  +0.12%  '
$ws.Range("D30").Value = "'1.86"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.69%  '
$ws.Range("D31").Value = "'25.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.15%  '
$ws.Range("D32").Value = "'1.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.36%  '
$ws.Range("D33").Value = "'2.25"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.37%  '
$ws.Range("D34").Value = "'5.36"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.15%  '
$ws.Range("D35").Value = "'54.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.32%  '
$ws.Range("D36").Value = "'5.76"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.98%  '
$ws.Range("D37").Value = "'439.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -11.01%  '
$ws.Range("D38").Value = '3.102.28'
$ws.Range("E38").Value = '  -4.43%  '
$ws.Range("D39").Value = "'0.0777"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.72%  '
$ws.Range("D40").Value = "'0.0372"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -7.72%  '
$ws.Range("D41").Value = "'0.115"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.43%  '
$ws.Range("D42").Value = "'7.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.32%  '
$ws.Range("D44").Value = "'2.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -13.21%  '
$ws.Range("D45").Value = "'25.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.31%  '
$ws.Range("D46").Value = "'0.236"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.36%  '
$ws.Range("D47").Value = "'0.106"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.88%  '
$ws.Range("B48").Value = 'BitgetToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb'
$ws.Range("D48").Value = "'1.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +9.84%  '
$ws.Range("B49").Value = 'Fetch.AI'
$ws.Range("C49").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D49").Value = "'1.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -7.22%  '
$ws.Range("D50").Value = "'113.35"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -8.97%  '
$ws.Range("D51").Value = '0.0₃0470'
$ws.Range("E51").Value = '  -11.17%  '
